$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: preserve the old last-row border styling (row 370) by copying it onto the new last row (390)
$ws.Range("A370:L370").Copy()
$ws.Range("A390").PasteSpecial(-4122)

# Step 2: row 370 is no longer the last row, so restyle it like a normal style-12-group row
$ws.Range("A366:L366").Copy()
$ws.Range("A370").PasteSpecial(-4122)

# Step 3: apply matching cell formats to each newly-inserted row (371-389)
$ws.Range("A369:L369").Copy()
$ws.Range("A371").PasteSpecial(-4122)
$ws.Range("A366:L366").Copy()
$ws.Range("A372").PasteSpecial(-4122)
$ws.Range("A367:I367").Copy()
$ws.Range("A373").PasteSpecial(-4122)
$ws.Range("M367:O367").Copy()
$ws.Range("M373").PasteSpecial(-4122)
$ws.Range("A366:L366").Copy()
$ws.Range("A374").PasteSpecial(-4122)
$ws.Range("A369:L369").Copy()
$ws.Range("A375").PasteSpecial(-4122)
$ws.Range("A368:I368").Copy()
$ws.Range("A376").PasteSpecial(-4122)
$ws.Range("M368:O368").Copy()
$ws.Range("M376").PasteSpecial(-4122)
$ws.Range("A369:L369").Copy()
$ws.Range("A377").PasteSpecial(-4122)
$ws.Range("A368:I368").Copy()
$ws.Range("A378").PasteSpecial(-4122)
$ws.Range("M368:O368").Copy()
$ws.Range("M378").PasteSpecial(-4122)
$ws.Range("A369:L369").Copy()
$ws.Range("A379").PasteSpecial(-4122)
$ws.Range("A366:L366").Copy()
$ws.Range("A380").PasteSpecial(-4122)
$ws.Range("A369:L369").Copy()
$ws.Range("A381").PasteSpecial(-4122)
$ws.Range("A366:L366").Copy()
$ws.Range("A382").PasteSpecial(-4122)
$ws.Range("A369:L369").Copy()
$ws.Range("A383").PasteSpecial(-4122)
$ws.Range("A366:L366").Copy()
$ws.Range("A384").PasteSpecial(-4122)
$ws.Range("A369:L369").Copy()
$ws.Range("A385").PasteSpecial(-4122)
$ws.Range("A368:I368").Copy()
$ws.Range("A386").PasteSpecial(-4122)
$ws.Range("M368:O368").Copy()
$ws.Range("M386").PasteSpecial(-4122)
$ws.Range("A369:L369").Copy()
$ws.Range("A387").PasteSpecial(-4122)
$ws.Range("A368:I368").Copy()
$ws.Range("A388").PasteSpecial(-4122)
$ws.Range("M368:O368").Copy()
$ws.Range("M388").PasteSpecial(-4122)
$ws.Range("A367:I367").Copy()
$ws.Range("A389").PasteSpecial(-4122)
$ws.Range("M367:O367").Copy()
$ws.Range("M389").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Step 4: write the response data for the new rows (371-390)
$ws.Cells.Item(371, 1).Value = 45594.407402129626
$ws.Cells.Item(371, 2).Value = "sbysooo@naver.com"
$ws.Cells.Item(371, 3).Value = "체육학과"
$ws.Cells.Item(371, 4).Value = 20227155
$ws.Cells.Item(371, 5).Value = "신윤수"
$ws.Cells.Item(371, 6).Value = "랜덤화"
$ws.Cells.Item(371, 7).Value = "28 vs 46"
$ws.Cells.Item(371, 8).Value = "플라시보 컨트롤 설계의 백신 접종 집단"
$ws.Cells.Item(371, 9).Value = "Red"
$ws.Cells.Item(371, 10).Value = "가. 10센트"
$ws.Cells.Item(371, 11).Value = "가. 100분"
$ws.Cells.Item(371, 12).Value = "가. 24일"

$ws.Cells.Item(372, 1).Value = 45594.40895590278
$ws.Cells.Item(372, 2).Value = "pjobin0821@naver.com"
$ws.Cells.Item(372, 3).Value = "러시아학과"
$ws.Cells.Item(372, 4).Value = 20241730
$ws.Cells.Item(372, 5).Value = "조유빈"
$ws.Cells.Item(372, 6).Value = "랜덤화"
$ws.Cells.Item(372, 7).Value = "28 vs 71"
$ws.Cells.Item(372, 8).Value = "NFIP 설계의 대조군 집단"
$ws.Cells.Item(372, 9).Value = "Red"
$ws.Cells.Item(372, 10).Value = "나. 5센트"
$ws.Cells.Item(372, 11).Value = "가. 100분"
$ws.Cells.Item(372, 12).Value = "나. 47일"

$ws.Cells.Item(373, 1).Value = 45594.45314092592
$ws.Cells.Item(373, 2).Value = "jiah1622@naver.com"
$ws.Cells.Item(373, 3).Value = "소프트웨어학부"
$ws.Cells.Item(373, 4).Value = 20245271
$ws.Cells.Item(373, 5).Value = "최지아"
$ws.Cells.Item(373, 6).Value = "가짜약 대조군"
$ws.Cells.Item(373, 7).Value = "28 vs 71"
$ws.Cells.Item(373, 8).Value = "플라시보 컨트롤 설계의 생리식염수 접종 집단"
$ws.Cells.Item(373, 9).Value = "Black"
$ws.Cells.Item(373, 13).Value = "나. 10센트"
$ws.Cells.Item(373, 14).Value = "나. 100분"
$ws.Cells.Item(373, 15).Value = "나. 24일"

$ws.Cells.Item(374, 1).Value = 45594.4776722338
$ws.Cells.Item(374, 2).Value = "r67890@naver.com"
$ws.Cells.Item(374, 3).Value = "스마트iot"
$ws.Cells.Item(374, 4).Value = 20205217
$ws.Cells.Item(374, 5).Value = "이규형"
$ws.Cells.Item(374, 6).Value = "랜덤화"
$ws.Cells.Item(374, 7).Value = "28 vs 25"
$ws.Cells.Item(374, 8).Value = "NFIP 설계의 백신 접종 집단"
$ws.Cells.Item(374, 9).Value = "Red"
$ws.Cells.Item(374, 10).Value = "가. 10센트"
$ws.Cells.Item(374, 11).Value = "나. 5분"
$ws.Cells.Item(374, 12).Value = "나. 47일"

$ws.Cells.Item(375, 1).Value = 45594.51212543981
$ws.Cells.Item(375, 2).Value = "jym85362@naver.com"
$ws.Cells.Item(375, 3).Value = "데이터사이언스"
$ws.Cells.Item(375, 4).Value = 20243234
$ws.Cells.Item(375, 5).Value = "유수현"
$ws.Cells.Item(375, 6).Value = "가짜약 대조군"
$ws.Cells.Item(375, 7).Value = "28 vs 25"
$ws.Cells.Item(375, 8).Value = "플라시보 컨트롤 설계의 생리식염수 접종 집단"
$ws.Cells.Item(375, 9).Value = "Red"
$ws.Cells.Item(375, 10).Value = "가. 10센트"
$ws.Cells.Item(375, 11).Value = "나. 5분"
$ws.Cells.Item(375, 12).Value = "가. 24일"

$ws.Cells.Item(376, 1).Value = 45594.53288690972
$ws.Cells.Item(376, 2).Value = "ksol902@naver.com"
$ws.Cells.Item(376, 3).Value = "데이터사이언스"
$ws.Cells.Item(376, 4).Value = 20243215
$ws.Cells.Item(376, 5).Value = "김은솔"
$ws.Cells.Item(376, 6).Value = "가짜약 대조군"
$ws.Cells.Item(376, 7).Value = "28 vs 46"
$ws.Cells.Item(376, 8).Value = "플라시보 컨트롤 설계의 백신 접종 집단"
$ws.Cells.Item(376, 9).Value = "Black"
$ws.Cells.Item(376, 13).Value = "가. 5센트"
$ws.Cells.Item(376, 14).Value = "가. 5분"
$ws.Cells.Item(376, 15).Value = "나. 24일"

$ws.Cells.Item(377, 1).Value = 45594.5554753125
$ws.Cells.Item(377, 2).Value = "jklucky09@naver.com"
$ws.Cells.Item(377, 3).Value = "사회학과 "
$ws.Cells.Item(377, 4).Value = 20171101
$ws.Cells.Item(377, 5).Value = "최준근 "
$ws.Cells.Item(377, 6).Value = "가짜약 대조군"
$ws.Cells.Item(377, 7).Value = "28 vs 71"
$ws.Cells.Item(377, 8).Value = "NFIP 설계의 대조군 집단"
$ws.Cells.Item(377, 9).Value = "Red"
$ws.Cells.Item(377, 10).Value = "가. 10센트"
$ws.Cells.Item(377, 11).Value = "나. 5분"
$ws.Cells.Item(377, 12).Value = "가. 24일"

$ws.Cells.Item(378, 1).Value = 45594.570895625
$ws.Cells.Item(378, 2).Value = "junseok5310@naver.com"
$ws.Cells.Item(378, 3).Value = "바이오메디컬학과"
$ws.Cells.Item(378, 4).Value = 20193626
$ws.Cells.Item(378, 5).Value = "오준석"
$ws.Cells.Item(378, 6).Value = "랜덤화"
$ws.Cells.Item(378, 7).Value = "28 vs 71"
$ws.Cells.Item(378, 8).Value = "NFIP 설계의 대조군 집단"
$ws.Cells.Item(378, 9).Value = "Black"
$ws.Cells.Item(378, 13).Value = "나. 10센트"
$ws.Cells.Item(378, 14).Value = "가. 5분"
$ws.Cells.Item(378, 15).Value = "가. 47일"

$ws.Cells.Item(379, 1).Value = 45594.717129247685
$ws.Cells.Item(379, 2).Value = "schoe357@gmail.com"
$ws.Cells.Item(379, 3).Value = "인공지능융합학부"
$ws.Cells.Item(379, 4).Value = 20246782
$ws.Cells.Item(379, 5).Value = "최성민"
$ws.Cells.Item(379, 6).Value = "랜덤화"
$ws.Cells.Item(379, 7).Value = "28 vs 71"
$ws.Cells.Item(379, 8).Value = "NFIP 설계의 대조군 집단"
$ws.Cells.Item(379, 9).Value = "Red"
$ws.Cells.Item(379, 10).Value = "가. 10센트"
$ws.Cells.Item(379, 11).Value = "가. 100분"
$ws.Cells.Item(379, 12).Value = "가. 24일"

$ws.Cells.Item(380, 1).Value = 45594.71713447917
$ws.Cells.Item(380, 2).Value = "dabinchoe05@gmail.com"
$ws.Cells.Item(380, 3).Value = "인공지능융합학부"
$ws.Cells.Item(380, 4).Value = 20246780
$ws.Cells.Item(380, 5).Value = "최다빈"
$ws.Cells.Item(380, 6).Value = "랜덤화"
$ws.Cells.Item(380, 7).Value = "28 vs 71"
$ws.Cells.Item(380, 8).Value = "NFIP 설계의 대조군 집단"
$ws.Cells.Item(380, 9).Value = "Red"
$ws.Cells.Item(380, 10).Value = "가. 10센트"
$ws.Cells.Item(380, 11).Value = "가. 100분"
$ws.Cells.Item(380, 12).Value = "가. 24일"

$ws.Cells.Item(381, 1).Value = 45594.73304571759
$ws.Cells.Item(381, 2).Value = "imhyeongu00@gmail.com"
$ws.Cells.Item(381, 3).Value = "언론방송융합미디어전공"
$ws.Cells.Item(381, 4).Value = 20192575
$ws.Cells.Item(381, 5).Value = "임현구"
$ws.Cells.Item(381, 6).Value = "랜덤화"
$ws.Cells.Item(381, 7).Value = "28 vs 71"
$ws.Cells.Item(381, 8).Value = "NFIP 설계의 대조군 집단"
$ws.Cells.Item(381, 9).Value = "Red"
$ws.Cells.Item(381, 10).Value = "가. 10센트"
$ws.Cells.Item(381, 11).Value = "나. 5분"
$ws.Cells.Item(381, 12).Value = "가. 24일"

$ws.Cells.Item(382, 1).Value = 45594.81128692129
$ws.Cells.Item(382, 2).Value = "jione0831@naver.com"
$ws.Cells.Item(382, 3).Value = "간호학과"
$ws.Cells.Item(382, 4).Value = 20246262
$ws.Cells.Item(382, 5).Value = "윤지원"
$ws.Cells.Item(382, 6).Value = "랜덤화"
$ws.Cells.Item(382, 7).Value = "28 vs 71"
$ws.Cells.Item(382, 8).Value = "NFIP 설계의 대조군 집단"
$ws.Cells.Item(382, 9).Value = "Red"
$ws.Cells.Item(382, 10).Value = "나. 5센트"
$ws.Cells.Item(382, 11).Value = "나. 5분"
$ws.Cells.Item(382, 12).Value = "나. 47일"

$ws.Cells.Item(383, 1).Value = 45594.871939652774
$ws.Cells.Item(383, 2).Value = "yeel6945@naver.com"
$ws.Cells.Item(383, 3).Value = "경영학과"
$ws.Cells.Item(383, 4).Value = 20203635
$ws.Cells.Item(383, 5).Value = "이수빈"
$ws.Cells.Item(383, 6).Value = "가짜약 대조군"
$ws.Cells.Item(383, 7).Value = "28 vs 71"
$ws.Cells.Item(383, 8).Value = "플라시보 컨트롤 설계의 백신 접종 집단"
$ws.Cells.Item(383, 9).Value = "Red"
$ws.Cells.Item(383, 10).Value = "가. 10센트"
$ws.Cells.Item(383, 11).Value = "나. 5분"
$ws.Cells.Item(383, 12).Value = "가. 24일"

$ws.Cells.Item(384, 1).Value = 45594.886753217594
$ws.Cells.Item(384, 2).Value = "joazzzzz@naver.com"
$ws.Cells.Item(384, 3).Value = "소프트웨어학부"
$ws.Cells.Item(384, 4).Value = 20203213
$ws.Cells.Item(384, 5).Value = "김진구"
$ws.Cells.Item(384, 6).Value = "랜덤화"
$ws.Cells.Item(384, 7).Value = "28 vs 71"
$ws.Cells.Item(384, 8).Value = "NFIP 설계의 대조군 집단"
$ws.Cells.Item(384, 9).Value = "Red"
$ws.Cells.Item(384, 10).Value = "나. 5센트"
$ws.Cells.Item(384, 11).Value = "나. 5분"
$ws.Cells.Item(384, 12).Value = "나. 47일"

$ws.Cells.Item(385, 1).Value = 45594.895077037036
$ws.Cells.Item(385, 2).Value = "yongwoo7701@gmail.com"
$ws.Cells.Item(385, 3).Value = "체육학과"
$ws.Cells.Item(385, 4).Value = 20244130
$ws.Cells.Item(385, 5).Value = "유용우"
$ws.Cells.Item(385, 6).Value = "랜덤화"
$ws.Cells.Item(385, 7).Value = "25 vs 54"
$ws.Cells.Item(385, 8).Value = "플라시보 컨트롤 설계의 백신 접종 집단"
$ws.Cells.Item(385, 9).Value = "Red"
$ws.Cells.Item(385, 10).Value = "나. 5센트"
$ws.Cells.Item(385, 11).Value = "가. 100분"
$ws.Cells.Item(385, 12).Value = "가. 24일"

$ws.Cells.Item(386, 1).Value = 45594.943638564815
$ws.Cells.Item(386, 2).Value = "cindy_lol@naver.com"
$ws.Cells.Item(386, 3).Value = "국어국문학전공"
$ws.Cells.Item(386, 4).Value = 20201007
$ws.Cells.Item(386, 5).Value = "권한별"
$ws.Cells.Item(386, 6).Value = "랜덤화"
$ws.Cells.Item(386, 7).Value = "28 vs 71"
$ws.Cells.Item(386, 8).Value = "NFIP 설계의 대조군 집단"
$ws.Cells.Item(386, 9).Value = "Black"
$ws.Cells.Item(386, 13).Value = "가. 5센트"
$ws.Cells.Item(386, 14).Value = "가. 5분"
$ws.Cells.Item(386, 15).Value = "가. 47일"

$ws.Cells.Item(387, 1).Value = 45594.949528171295
$ws.Cells.Item(387, 2).Value = "minsung5342@naver.com"
$ws.Cells.Item(387, 3).Value = "사회복지학과"
$ws.Cells.Item(387, 4).Value = 20232311
$ws.Cells.Item(387, 5).Value = "김민성"
$ws.Cells.Item(387, 6).Value = "랜덤화"
$ws.Cells.Item(387, 7).Value = "28 vs 71"
$ws.Cells.Item(387, 8).Value = "NFIP 설계의 대조군 집단"
$ws.Cells.Item(387, 9).Value = "Red"
$ws.Cells.Item(387, 10).Value = "나. 5센트"
$ws.Cells.Item(387, 11).Value = "나. 5분"
$ws.Cells.Item(387, 12).Value = "나. 47일"

$ws.Cells.Item(388, 1).Value = 45594.95650006944
$ws.Cells.Item(388, 2).Value = "ohsolbi050521@gmail.com"
$ws.Cells.Item(388, 3).Value = "미래융합스쿨"
$ws.Cells.Item(388, 4).Value = 20246631
$ws.Cells.Item(388, 5).Value = "오솔비"
$ws.Cells.Item(388, 6).Value = "랜덤화"
$ws.Cells.Item(388, 7).Value = "28 vs 71"
$ws.Cells.Item(388, 8).Value = "NFIP 설계의 대조군 집단"
$ws.Cells.Item(388, 9).Value = "Black"
$ws.Cells.Item(388, 13).Value = "가. 5센트"
$ws.Cells.Item(388, 14).Value = "가. 5분"
$ws.Cells.Item(388, 15).Value = "가. 47일"

$ws.Cells.Item(389, 1).Value = 45595.07543189815
$ws.Cells.Item(389, 2).Value = "lucas3767@naver.com"
$ws.Cells.Item(389, 3).Value = "법학과"
$ws.Cells.Item(389, 4).Value = 20242702
$ws.Cells.Item(389, 5).Value = "권민재"
$ws.Cells.Item(389, 6).Value = "랜덤화"
$ws.Cells.Item(389, 7).Value = "28 vs 71"
$ws.Cells.Item(389, 8).Value = "NFIP 설계의 대조군 집단"
$ws.Cells.Item(389, 9).Value = "Black"
$ws.Cells.Item(389, 13).Value = "나. 10센트"
$ws.Cells.Item(389, 14).Value = "나. 100분"
$ws.Cells.Item(389, 15).Value = "나. 24일"

$ws.Cells.Item(390, 1).Value = 45595.1459509375
$ws.Cells.Item(390, 2).Value = "msy123581@gmail.com"
$ws.Cells.Item(390, 3).Value = "중국학과"
$ws.Cells.Item(390, 4).Value = 20241519
$ws.Cells.Item(390, 5).Value = "문신영"
$ws.Cells.Item(390, 6).Value = "랜덤화"
$ws.Cells.Item(390, 7).Value = "28 vs 71"
$ws.Cells.Item(390, 8).Value = "플라시보 컨트롤 설계의 생리식염수 접종 집단"
$ws.Cells.Item(390, 9).Value = "Red"
$ws.Cells.Item(390, 10).Value = "나. 5센트"
$ws.Cells.Item(390, 11).Value = "나. 5분"
$ws.Cells.Item(390, 12).Value = "나. 47일"

# Step 5: grow the table (ListObject) to cover the newly-added rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:O390"))
